# Update the "timestamp" column (O) for all data rows (2-73) from the
# old scrape time to the new scrape time, matching the commit
# "Upload excel files with prices".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 73; $row++) {
    $cell = $ws.Range("O$row")
    if ($cell.Text -eq "2022-08-11 07:01:18") {
        $cell.Value = "2022-08-11 20:57:09"
    }
}
